$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("R11").Value = 328
$ws.Range("S11").Value = 36.419
$ws.Range("T11").Value = 1635
$ws.Range("U11").Value = 181.538

# Row 18
$ws.Range("R18").Value = 367
$ws.Range("S18").Value = 31.666
$ws.Range("T18").Value = 1892
$ws.Range("U18").Value = 163.25
$ws.Range("V18").Value = 2533.312
$ws.Range("W18").Value = 218.585
$ws.Range("X18").Value = 814.423
$ws.Range("Y18").Value = 70.272

# Row 27
$ws.Range("R27").Value = 340
$ws.Range("S27").Value = 48.932
$ws.Range("T27").Value = 3485
$ws.Range("U27").Value = 501.551

# Row 32
$ws.Range("R32").Value = 859
$ws.Range("S32").Value = 22.76
$ws.Range("T32").Value = 4490
$ws.Range("U32").Value = 118.965

# Row 43
$ws.Range("T43").Value = 1761
$ws.Range("U43").Value = 428.961
$ws.Range("X43").Value = 978.92
$ws.Range("Y43").Value = 238.455

# Row 45
$ws.Range("R45").Value = 39
$ws.Range("S45").Value = 44.526
$ws.Range("T45").Value = 208
$ws.Range("U45").Value = 237.47
$ws.Range("V45").Value = 17
$ws.Range("W45").Value = 19.409
$ws.Range("X45").Value = 112
$ws.Range("Y45").Value = 127.869

# Row 46
$ws.Range("R46").Value = 1119
$ws.Range("S46").Value = 104.492
$ws.Range("T46").Value = 6449
$ws.Range("U46").Value = 602.205
$ws.Range("V46").Value = 1985.975
$ws.Range("W46").Value = 185.449
$ws.Range("X46").Value = 11397.99
$ws.Range("Y46").Value = 1064.339

# Row 48
$ws.Range("R48").Value = 131
$ws.Range("S48").Value = 22.617
$ws.Range("T48").Value = 799
$ws.Range("U48").Value = 137.944
$ws.Range("X48").Value = 667.401
$ws.Range("Y48").Value = 115.224

# Row 57
$ws.Range("R57").Value = 45
$ws.Range("S57").Value = 33.923
$ws.Range("T57").Value = 433
$ws.Range("U57").Value = 326.413
$ws.Range("V57").Value = 29.038
$ws.Range("W57").Value = 21.89
$ws.Range("X57").Value = 298.387
$ws.Range("Y57").Value = 224.936

# Row 61
$ws.Range("R61").Value = 28
$ws.Range("S61").Value = 5.053
$ws.Range("T61").Value = 145
$ws.Range("U61").Value = 26.17

# Row 62
$ws.Range("R62").Value = 2766
$ws.Range("S62").Value = 42.376
$ws.Range("T62").Value = 25235
$ws.Range("U62").Value = 386.604
$ws.Range("V62").Value = 1354.896
$ws.Range("W62").Value = 20.757
$ws.Range("X62").Value = 9104.391
$ws.Range("Y62").Value = 139.481

# Row 66
$ws.Range("R66").Value = 4950
$ws.Range("S66").Value = 59.081
$ws.Range("X66").Value = 2333.297
$ws.Range("Y66").Value = 27.849

# Row 68
$ws.Range("V68").Value = 102.048
$ws.Range("W68").Value = 9.791

# Row 77
$ws.Range("T77").Value = 4345
$ws.Range("U77").Value = 449.777

# Row 78
$ws.Range("T78").Value = 19
$ws.Range("U78").Value = 55.678
$ws.Range("X78").Value = 0.956
$ws.Range("Y78").Value = 2.801

# Row 84
$ws.Range("R84").Value = 196
$ws.Range("S84").Value = 39.694
$ws.Range("T84").Value = 1923
$ws.Range("U84").Value = 389.445
$ws.Range("V84").Value = 59.404
$ws.Range("W84").Value = 12.03
$ws.Range("X84").Value = 900.117
$ws.Range("Y84").Value = 182.291

# Row 86
$ws.Range("R86").Value = 2503
$ws.Range("S86").Value = 41.398
$ws.Range("T86").Value = 25260
$ws.Range("U86").Value = 417.784
$ws.Range("X86").Value = 3871.549
$ws.Range("Y86").Value = 64.033

# Row 96
$ws.Range("T96").Value = 1136
$ws.Range("U96").Value = 602.268
$ws.Range("V96").Value = 148.344
$ws.Range("W96").Value = 78.647
$ws.Range("X96").Value = 794.772
$ws.Range("Y96").Value = 421.361

# Row 102
$ws.Range("T102").Value = 1910
$ws.Range("U102").Value = 701.615
$ws.Range("V102").Value = 95.479
$ws.Range("W102").Value = 35.073
$ws.Range("X102").Value = 514.415
$ws.Range("Y102").Value = 188.964

# Row 103
$ws.Range("R103").Value = 20
$ws.Range("S103").Value = 31.95
$ws.Range("T103").Value = 88
$ws.Range("U103").Value = 140.58

# Row 109
$ws.Range("X109").Value = 42.941
$ws.Range("Y109").Value = 97.253

# Row 124
$ws.Range("R124").Value = 676
$ws.Range("S124").Value = 39.452
$ws.Range("T124").Value = 1708
$ws.Range("U124").Value = 99.68
$ws.Range("V124").Value = 236.963
$ws.Range("W124").Value = 13.829
$ws.Range("X124").Value = 1336.511
$ws.Range("Y124").Value = 77.999

# Row 130
$ws.Range("T130").Value = 158
$ws.Range("U130").Value = 29.145
$ws.Range("V130").Value = 14.244
$ws.Range("W130").Value = 2.628
$ws.Range("X130").Value = 91.571
$ws.Range("Y130").Value = 16.891

# Row 139
$ws.Range("T139").Value = 15588
$ws.Range("U139").Value = 411.873

# Row 140
$ws.Range("R140").Value = 647
$ws.Range("S140").Value = 63.452
$ws.Range("T140").Value = 4889
$ws.Range("U140").Value = 479.469
$ws.Range("X140").Value = 549.692
$ws.Range("Y140").Value = 53.909

# Row 142
$ws.Range("R142").Value = 1076
$ws.Range("S142").Value = 55.932
$ws.Range("X142").Value = 8379.005
$ws.Range("Y142").Value = 435.552

# Row 157
$ws.Range("T157").Value = 3164
$ws.Range("U157").Value = 579.525

# Row 158
$ws.Range("R158").Value = 191
$ws.Range("S158").Value = 91.874
$ws.Range("T158").Value = 1220
$ws.Range("U158").Value = 586.84
$ws.Range("V158").Value = 106.898
$ws.Range("W158").Value = 51.42
$ws.Range("X158").Value = 747.29
$ws.Range("Y158").Value = 359.458

# Row 164
$ws.Range("R164").Value = 2651
$ws.Range("S164").Value = 56.7
$ws.Range("T164").Value = 17645
$ws.Range("U164").Value = 377.395
$ws.Range("V164").Value = 136.468
$ws.Range("W164").Value = 2.919
$ws.Range("X164").Value = 3078.996
$ws.Range("Y164").Value = 65.854

# Row 168
$ws.Range("R168").Value = 354
$ws.Range("S168").Value = 35.052
$ws.Range("T168").Value = 2531
$ws.Range("U168").Value = 250.612
$ws.Range("V168").Value = 139.196
$ws.Range("W168").Value = 13.783

# Row 183
$ws.Range("R183").Value = 3953
$ws.Range("S183").Value = 58.23
$ws.Range("T183").Value = 38676
$ws.Range("U183").Value = 569.72
$ws.Range("X183").Value = 28362
$ws.Range("Y183").Value = 417.789

# Row 184
$ws.Range("R184").Value = 22304
$ws.Range("S184").Value = 67.383
$ws.Range("T184").Value = 119927
$ws.Range("U184").Value = 362.314
$ws.Range("V184").Value = 1158
$ws.Range("W184").Value = 3.498
$ws.Range("X184").Value = 27012
$ws.Range("Y184").Value = 81.607
